$wb = $excel.ActiveWorkbook

# Sheet 1: two_blank_rows
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A3").Value = "var1"
$ws1.Range("B3").Value = "var2"
$ws1.Range("A4:B4").ClearContents()
$ws1.Range("A5").Value = "v2,1"
$ws1.Range("B5").Value = "v2,2"
$ws1.Range("A6:B6").ClearContents()
$ws1.Range("A7").Value = "v4,1"
$ws1.Range("B7").Value = "v4,2"

# Sheet 2: occupied_row_and_blank_row
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A3").Value = "var1"
$ws2.Range("B3").Value = "var2"
$ws2.Range("A4:B4").ClearContents()
$ws2.Range("A5").Value = "v2,1"
$ws2.Range("B5").Value = "v2,2"
$ws2.Range("A6:B6").ClearContents()
$ws2.Range("A7").Value = "v4,1"
$ws2.Range("B7").Value = "v4,2"

# Sheet 3: two_occupied_rows
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A3").Value = "var1"
$ws3.Range("B3").Value = "var2"
$ws3.Range("A4:B4").ClearContents()
$ws3.Range("A5").Value = "v2,1"
$ws3.Range("B5").Value = "v2,2"
$ws3.Range("A6:B6").ClearContents()
$ws3.Range("A7").Value = "v4,1"
$ws3.Range("B7").Value = "v4,2"

# Update selections on each sheet; select sheet3 last so it remains the
# active tab (matching tabSelected="1" staying on two_occupied_rows).
$ws1.Range("A3:B7").Select()
$ws2.Range("A3:B7").Select()
$ws3.Range("A3:B7").Select()
